$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to be treated as literal text (so date-like and
    # numeric-looking strings such as "2024-04-07" or "235" are not
    # reinterpreted as a date/number), then restore the default "Normal"
    # style so no stray number-format style sticks to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 (Sujan) - update date, purpose, amount
Set-TextValue $ws.Range("B2") "2024-04-07"
$ws.Range("C2").Value = "Dinner"
Set-TextValue $ws.Range("D2") "235"

# Row 3 (Hari) - update date, purpose, amount
Set-TextValue $ws.Range("B3") "2024-04-09"
$ws.Range("C3").Value = "Breakfast"
Set-TextValue $ws.Range("D3") "340"

# Row 4 - name change (Ram Chaudhary -> Asmin Dhakal), date, amount
$ws.Range("A4").Value = "Asmin Dhakal"
Set-TextValue $ws.Range("B4") "2024-04-01"
Set-TextValue $ws.Range("D4") "5000"

# Row 5 (new) - admin
$ws.Range("A5").Value = "admin"
Set-TextValue $ws.Range("B5") "2024-03-31"
$ws.Range("C5").Value = "Breakfast"
Set-TextValue $ws.Range("D5") "245"

# Row 6 (new) - Asmin Dhakal
$ws.Range("A6").Value = "Asmin Dhakal"
Set-TextValue $ws.Range("B6") "2024-04-08"
$ws.Range("C6").Value = "Breakfast"
Set-TextValue $ws.Range("D6") "200"
